# Update NATMI LR-pair worksheet with recalculated TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lta", "Tnfrsf1b", "ECs", 1, 0.3333333333333333, 0.07057933333333334, 0.211738, 0.1281663575568867, 0.1281663575568867, 3, 1, 11.53121, 34.59363, 0.04670311854310271, 0.04697417460197403, 0.8138651143266667, 7.32478602894, 0.005985768590216968, 0.006020508857976231),
    @("ECs", "Lta", "Tnfrsf1b", "FAPs", 1, 0.3333333333333333, 0.07057933333333334, 0.211738, 0.1281663575568867, 0.1281663575568867, 3, 1, 12.200174, 36.600522, 0.04941252241252041, 0.04969930333854504, 0.8610801474706666, 7.749721327236, 0.00633302301531077, 0.006369778682016137),
    @("ECs", "Lta", "Tnfrsf1b", "Inflammatory-Mac", 1, 0.3333333333333333, 0.07057933333333334, 0.211738, 0.1281663575568867, 0.1281663575568867, 3, 1, 113.969907, 341.909721, 0.4615951038067463, 0.4642741144067919, 8.043920056121999, 72.395280505098, 0.05916096312100368, 0.05950432215146782),
    @("ECs", "Lta", "Tnfrsf1b", "MuSCs", 1, 0.3333333333333333, 0.07057933333333334, 0.211738, 0.1281663575568867, 0.1281663575568867, 2, 1, 4.2741545, 8.548309, 0.01731096253429049, 0.01160762138918714, 0.3016669751736667, 1.810001851042, 0.002218683013823745, 0.001487706553351526),
    @("ECs", "Lta", "Tnfrsf1b", "Resolving-Mac", 1, 0.3333333333333333, 0.07057933333333334, 0.211738, 0.1281663575568867, 0.1281663575568867, 3, 1, 104.9290516666667, 314.787155, 0.4249782927033401, 0.4274447862635018, 7.405822513932223, 66.65240262539, 0.05446791981653155, 0.05478404131207499),
    @("FAPs", "Lta", "Tnfrsf1b", "ECs", 1, 0.3333333333333333, 0.0279, 0.0837, 0.05066414213561767, 0.05066414213561768, 3, 1, 11.53121, 34.59363, 0.04670311854310271, 0.04697417460197403, 0.321720759, 2.895486830999999, 0.002366173436044357, 0.002379906258737735),
    @("FAPs", "Lta", "Tnfrsf1b", "FAPs", 1, 0.3333333333333333, 0.0279, 0.0837, 0.05066414213561767, 0.05066414213561768, 3, 1, 12.200174, 36.600522, 0.04941252241252041, 0.04969930333854504, 0.3403848546, 3.0634636914, 0.002503443058787328, 0.002517972568385224),
    @("FAPs", "Lta", "Tnfrsf1b", "Inflammatory-Mac", 1, 0.3333333333333333, 0.0279, 0.0837, 0.05066414213561767, 0.05066414213561768, 3, 1, 113.969907, 341.909721, 0.4615951038067463, 0.4642741144067919, 3.1797604053, 28.6178436477, 0.02338631994837019, 0.02352204972219373),
    @("FAPs", "Lta", "Tnfrsf1b", "MuSCs", 1, 0.3333333333333333, 0.0279, 0.0837, 0.05066414213561767, 0.05066414213561768, 2, 1, 4.2741545, 8.548309, 0.01731096253429049, 0.01160762138918714, 0.11924891055, 0.7154934633, 0.0008770450663416458, 0.0005880901799182135),
    @("FAPs", "Lta", "Tnfrsf1b", "Resolving-Mac", 1, 0.3333333333333333, 0.0279, 0.0837, 0.05066414213561767, 0.05066414213561768, 3, 1, 104.9290516666667, 314.787155, 0.4249782927033401, 0.4274447862635018, 2.9275205415, 26.3476848735, 0.02153116062607415, 0.02165612340638277),
    @("Inflammatory-Mac", "Lta", "Tnfrsf1b", "ECs", 1, 0.3333333333333333, 0.128589, 0.385767, 0.2335072176730087, 0.2335072176730087, 3, 1, 11.53121, 34.59363, 0.04670311854310271, 0.04697417460197403, 1.48278676269, 13.34508086421, 0.01090551526765261, 0.01096880881379307),
    @("Inflammatory-Mac", "Lta", "Tnfrsf1b", "FAPs", 1, 0.3333333333333333, 0.128589, 0.385767, 0.2335072176730087, 0.2335072176730087, 3, 1, 12.200174, 36.600522, 0.04941252241252041, 0.04969930333854504, 1.568808174486, 14.119273570374, 0.01153818062675282, 0.01160514604287053),
    @("Inflammatory-Mac", "Lta", "Tnfrsf1b", "Inflammatory-Mac", 1, 0.3333333333333333, 0.128589, 0.385767, 0.2335072176730087, 0.2335072176730087, 3, 1, 113.969907, 341.909721, 0.4615951038067463, 0.4642741144067919, 14.655276371223, 131.897487341007, 0.1077857883813969, 0.1084113566927301),
    @("Inflammatory-Mac", "Lta", "Tnfrsf1b", "MuSCs", 1, 0.3333333333333333, 0.128589, 0.385767, 0.2335072176730087, 0.2335072176730087, 2, 1, 4.2741545, 8.548309, 0.01731096253429049, 0.01160762138918714, 0.5496092530005, 3.297655518003, 0.004042234696623868, 0.002710463374390794),
    @("Inflammatory-Mac", "Lta", "Tnfrsf1b", "Resolving-Mac", 1, 0.3333333333333333, 0.128589, 0.385767, 0.2335072176730087, 0.2335072176730087, 3, 1, 104.9290516666667, 314.787155, 0.4249782927033401, 0.4274447862635018, 13.492721824765, 121.434496422885, 0.09923549870058243, 0.0998114427492242),
    @("MuSCs", "Lta", "Tnfrsf1b", "ECs", 3, 1, 0.323617, 0.970851, 0.5876622826344869, 0.5876622826344869, 3, 1, 11.53121, 34.59363, 0.04670311854310271, 0.04697417460197403, 3.73169558657, 33.58526027913, 0.02744566124918877, 0.027604950671467),
    @("MuSCs", "Lta", "Tnfrsf1b", "FAPs", 3, 1, 0.323617, 0.970851, 0.5876622826344869, 0.5876622826344869, 3, 1, 12.200174, 36.600522, 0.04941252241252041, 0.04969930333854504, 3.948183709357999, 35.533653384222, 0.02903787571166949, 0.02920640604527315),
    @("MuSCs", "Lta", "Tnfrsf1b", "Inflammatory-Mac", 3, 1, 0.323617, 0.970851, 0.5876622826344869, 0.5876622826344869, 3, 1, 113.969907, 341.909721, 0.4615951038067463, 0.4642741144067919, 36.882599393619, 331.943394542571, 0.2712620323559755, 0.2728363858404003),
    @("MuSCs", "Lta", "Tnfrsf1b", "MuSCs", 3, 1, 0.323617, 0.970851, 0.5876622826344869, 0.5876622826344869, 2, 1, 4.2741545, 8.548309, 0.01731096253429049, 0.01160762138918714, 1.3831890568265, 8.299134340959, 0.01017299975750123, 0.006821361281526612),
    @("MuSCs", "Lta", "Tnfrsf1b", "Resolving-Mac", 3, 1, 0.323617, 0.970851, 0.5876622826344869, 0.5876622826344869, 3, 1, 104.9290516666667, 314.787155, 0.4249782927033401, 0.4274447862635018, 33.95682491321166, 305.611424218905, 0.249743713560152, 0.2511931787958199)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 2 + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
